$d = $word.ActiveDocument

function Get-ParaIndexForRange($rng) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        if ($pp.Range.Start -le $rng.Start -and $pp.Range.End -ge $rng.End) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Insert a new bullet "Duplicated records (had 8 pairs of complete
#    duplicates)" right before the "Records with surge percent..." bullet.
# ------------------------------------------------------------------
$rngAnchor = $d.Content.Duplicate
$null = $rngAnchor.Find.Execute("Records with surge percent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIdx = Get-ParaIndexForRange $rngAnchor
$anchorPara = $d.Paragraphs.Item($anchorIdx)

# InsertParagraphBefore creates a new empty paragraph immediately before
# this one, inheriting its paragraph formatting (style / list / spacing).
$anchorPara.Range.InsertParagraphBefore()

# The freshly created paragraph now occupies the anchor's old slot.
$newBulletPara = $d.Paragraphs.Item($anchorIdx)
$newBulletPara.Range.Text = "Duplicated records (had 8 pairs of complete duplicates)"

# ------------------------------------------------------------------
# 2) Split "dropped 2,133 records leaving 96% of the original data."
#    into three runs: "dropped 2,1" / "41" / " records leaving 96% of
#    the original data." (i.e. the figure 2,133 becomes 2,141).
# ------------------------------------------------------------------
$rngOld = $d.Content.Duplicate
$null = $rngOld.Find.Execute("dropped 2,133 records leaving 96% of the original data.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hostIdx = Get-ParaIndexForRange $rngOld
$hostPara = $d.Paragraphs.Item($hostIdx)

# Remove the old single run's text entirely, leaving "These changes ".
$rngOld.Text = ""

# Split the (now shortened) host paragraph so the replacement text can be
# inserted into a brand new paragraph - this lets us use InsertXML to
# build exact run boundaries without the engine silently re-coalescing
# runs that live in a paragraph that also received a plain-text edit.
$hostPara.Range.InsertParagraphAfter()
$tailPara = $d.Paragraphs.Item($hostIdx + 1)

$openXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:t>dropped 2,1</w:t></w:r><w:r><w:t>41</w:t></w:r><w:r><w:t xml:space="preserve"> records leaving 96% of the original data.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailPara.Range.InsertXML($openXml)

# Merge the host paragraph and the tail paragraph back into a single
# paragraph by deleting the paragraph mark between them - the tail
# paragraph (which now carries the correct pPr + the three new runs)
# survives, giving us exactly: "These changes " + 3 split runs.
$markRange = $d.Range($hostPara.Range.End - 1, $hostPara.Range.End)
$markRange.Delete()
